$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (and workbook.xml <sheet name=.../>)
$ws.Name = "ShearF"

# Copy the formatting (bold font, borders, centered alignment) from A15 onto
# the new row's A cell before filling in values, matching style index 1
# used throughout column A / header rows.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

# Populate the new row (row 16) of data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9768042872944713
$ws.Range("D16").Value = 1.096363831534633
$ws.Range("E16").Value = 0.9698176992801386
$ws.Range("F16").Value = 0.9768042872944713
$ws.Range("G16").Value = 1.051940079293056
$ws.Range("H16").Value = 0.9281321562220225
$ws.Range("I16").Value = 0.9681114896547558
$ws.Range("J16").Value = 1.096363831534633
$ws.Range("K16").Value = 1.033090765407386
$ws.Range("L16").Value = 1.004947526350928
$ws.Range("M16").Value = 0.9985282572131796
